# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the first row (08baa50c-b7f6-462f-9e13-1e2434a3f926.md) across the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 03:10:00"

# zh-cn sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K) for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 03:09:55"
$wsZhCn.Range("K2").Value = "2016-09-04 03:10:19"

# de-de sheet: Correspond Handoff Datetime (H) / Correspond Handback DateTime (K) for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 03:10:00"
$wsDeDe.Range("K2").Value = "2016-09-04 03:10:26"
